# Scraper run @ 31/12/2025 07:57:39 — appends new arrival rows to each of
# the three sheets (LP1912, LP1912-215, 6203-6173) and refreshes the
# "Última actualización" / "Total filas" header cells on every sheet.

$wb = $excel.ActiveWorkbook

$timestamp = "Última actualización: 31/12/2025 07:57:39"

# ---------------------------------------------------------------------
# Sheet "LP1912": columns A(-), B=Hora_Scrap, C=Hora_Llegada, D=Línea,
# E=Minutos, F=Parada, G=Fecha
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = $timestamp
$ws1.Range("A3").Value = "Total filas: 742"

$rows1 = @(
    @("08:01", "16_SANTA ANA", 4),
    @("08:03", "17_ROMERO", 6),
    @("08:11", "16_SANTA ANA", 14),
    @("08:14", "10_OLMOS", 17),
    @("08:29", "14_ABASTO", 32),
    @("08:44", "10_OLMOS", 47),
    @("08:51", "16_SANTA ANA", 54),
    @("09:02", "17X38_ROMERO", 65),
    @("09:03", "23_HERNANDEZ", 66),
    @("09:14", "11_ETCHEVERRY", 77),
    @("09:17", "27_EL RETIRO", 80),
    @("09:27", "215_EL PELIGRO", 90),
    @("09:33", "23_HERNANDEZ", 96)
)

$r = 731
foreach ($row in $rows1) {
    $ws1.Cells.Item($r, 2).Value = "07:57:29"
    $ws1.Cells.Item($r, 3).Value = $row[0]
    $ws1.Cells.Item($r, 4).Value = $row[1]
    $ws1.Cells.Item($r, 5).Value = $row[2]
    $ws1.Cells.Item($r, 6).Value = "LP1912"
    $ws1.Cells.Item($r, 7).Value = "31/12/2025"
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet "LP1912-215": columns A(-), B=Fecha, C=Hora_Scrap, D=Hora_Llegada,
# E=Línea, F=Minutos, G=Parada
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = $timestamp
$ws2.Range("A3").Value = "Total filas: 53"

$ws2.Cells.Item(54, 2).Value = "31/12/2025"
$ws2.Cells.Item(54, 3).Value = "07:57:29"
$ws2.Cells.Item(54, 4).Value = "09:27"
$ws2.Cells.Item(54, 5).Value = "215_EL PELIGRO"
$ws2.Cells.Item(54, 6).Value = 90
$ws2.Cells.Item(54, 7).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "6203-6173": columns A(-), B=Fecha, C=Hora_Scrap, D=Hora_Llegada,
# E=Línea, F=Minutos, G=Parada
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = $timestamp
$ws3.Range("A3").Value = "Total filas: 90"

$ws3.Cells.Item(91, 2).Value = "31/12/2025"
$ws3.Cells.Item(91, 3).Value = "07:57:39"
$ws3.Cells.Item(91, 4).Value = "08:10"
$ws3.Cells.Item(91, 5).Value = "215A_LA PLATA"
$ws3.Cells.Item(91, 6).Value = 13
$ws3.Cells.Item(91, 7).Value = "L6173"
